$wb = $excel.ActiveWorkbook

# The "EN" sheet lists club names; two clubs (Millwall at row 13 and QPR at
# row 16) are removed, shifting the remaining rows up.
$ws = $wb.Worksheets.Item("EN")

# Delete row 13 (Millwall) first; QPR (originally row 16) is now at row 15.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(15).Delete()

# EN becomes the active/selected sheet & tab, with B15 selected.
$ws.Range("B15").Select() | Out-Null
$ws.Activate() | Out-Null
